$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.732581377029419
$ws.Range("B1").Value = 4.256258010864258
$ws.Range("C1").Value = 3.381970167160034
$ws.Range("D1").Value = 2.181460618972778
$ws.Range("E1").Value = 1.882538080215454
